# SORMAS User Rights workbook update
# - Adds 4 new user rights to the "User Rights" sheet (rows 74-77):
#     CLINICAL_COURSE_VIEW, CLINICAL_VISIT_CREATE, CLINICAL_VISIT_EDIT, CLINICAL_VISIT_DELETE
# - Bumps the SORMAS Version string on the "About" sheet to 1.13.0-SNAPSHOT

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$about = $wb.Worksheets.Item(2)

$yes = "Yes"
$no = "No"

# Column order on the sheet: A User Right, B Description, then one Yes/No column per role:
# C Admin, D National User, E Surveillance Supervisor, F Surveillance Officer,
# G Hospital Informant, H Community Informant, I Clinician, J Case Officer,
# K Contact Supervisor, L Contact Officer, M Event Officer, N Lab Officer,
# O External Lab Officer, P National Observer, Q State Observer, R District Observer
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

# New rows to append, keyed by user right name -> per-role Yes/No pattern
$newRights = @(
    @{ Name = "CLINICAL_COURSE_VIEW";  Pattern = @(1,1,0,0,0,0,1,0,0,0,0,0,0,1,1,1) },
    @{ Name = "CLINICAL_VISIT_CREATE"; Pattern = @(1,1,0,0,0,0,1,0,0,0,0,0,0,0,0,0) },
    @{ Name = "CLINICAL_VISIT_EDIT";   Pattern = @(1,1,0,0,0,0,1,0,0,0,0,0,0,0,0,0) },
    @{ Name = "CLINICAL_VISIT_DELETE"; Pattern = @(1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0) }
)

# Insert the new rows right after the last existing data row (73), copying that
# row's formatting (fills/fonts/borders) down so the new rows look consistent
# with the rest of the table.
$insertAt = 74
$lastRow = 73
for ($i = 0; $i -lt $newRights.Count; $i++) {
    $ws.Range("A$lastRow`:R$lastRow").Copy() | Out-Null
    $ws.Range("A$insertAt`:R$insertAt").Insert() | Out-Null
}

# Re-apply a clean thin black border to the pasted-in rows (the plain row
# insert above does not carry the border through), matching the look of the
# rest of the table.
$borderRng = $ws.Range("C$insertAt`:R" + ($insertAt + $newRights.Count - 1))
$borderRng.Borders.Color = 0
$borderRng.Borders.LineStyle = 1

for ($i = 0; $i -lt $newRights.Count; $i++) {
    $r = $insertAt + $i
    $right = $newRights[$i]

    $ws.Range("A$r").Value = $right.Name
    $ws.Range("B$r").Value = $right.Name

    for ($c = 0; $c -lt $cols.Count; $c++) {
        $cell = $ws.Range($cols[$c] + $r)
        if ($right.Pattern[$c] -eq 1) {
            $cell.Value = $yes
        } else {
            $cell.Value = $no
        }
    }
}

# Bump the SORMAS Version shown on the About sheet
$about.Range("A2").Value = "1.13.0-SNAPSHOT"
